$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time log entry: 09/14/2023, Internship, Completed 8 hours assisting with daily operations
$ws.Range("A4").Value = 45183
$ws.Range("A4").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("B4").Value = "Internship"
$ws.Range("C4").Value = "Completed 8 hours assisting with daily operations"

$ws.Range("D4").Select()
